$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C for data rows 2-201 advances
# by one day, from serial 45205 (2023-10-06) to 45206 (2023-10-07).
$ws.Range("C2:C201").Value = 45206
